$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Double Bass transpose value (H18: 0 -> 12) ---
$ws.Range("H18").Value = 12

# --- Add new row 19: Electric Bass ---
# Values are set in an order that mirrors how the underlying shared-string
# table was built (en/de names first, the normalised "key" last).
$ws.Range("B19").Value = "Electric Bass"
$ws.Range("C19").Value = "NA"
$ws.Range("D19").Value = "E-Bass"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = 12
$ws.Range("I19").Value = "bass"
$ws.Range("J19").Value = 40
$ws.Range("A19").Value = "electric_bass"

# --- Normalise formatting so the whole "de"/"lv" column (E) and the left
# "key" column (A, rows 14-19) as well as the new J19 cell use the same
# Arial 12 black style already used elsewhere in the sheet (D3). Using
# copy / paste-special (formats only) re-uses an existing style entry
# instead of fabricating new ones. ---
$ws.Range("D3").Copy()
$ws.Range("E3:E14").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("A14:A19").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("J19").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Make the header row bold ---
$ws.Range("A1:J1").Font.Bold = $true

# --- Update sheet view / selection state ---
$ws.Range("A19").Select() | Out-Null
